$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date column (C) for rows 2-11 from 2023-09-20 (45189) to 2023-09-21 (45190)
for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45189) {
        $cell.Value = 45190
    }
}
